$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers
$ws.Range("Q2").Value = 522980
$ws.Range("R2").Value = 7077151

# Remove the Starttid (Z2) and Sluttid (AB2) cell contents entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
